$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "grant_en"

# Clear existing content first
$ws.Range("A1:E30").ClearContents()

# Write new cell values
$ws.Range("A1").Value = 'what'
$ws.Range("B1").Value = 'when'
$ws.Range("C1").Value = 'with'
$ws.Range("D1").Value = 'where'
$ws.Range("E1").Value = 'why'

$ws.Range("A2").Value = '\href{https://minciencias.gov.co/convocatorias/construccion-paz-programa-y-proyectos-ctei-fortalecimiento-capacidades-para-la}{Postdoctoral Research Stays -  Call 935-2023 - Orchids Program. Women in science: agents for peace: Agents for Peace 2023}'
$ws.Range("B2").Value = 'Dic. 2023 - Jan. 2025'
$ws.Range("C2").Value = '\href{https://minciencias.gov.co/}{Minciencias}'
$ws.Range("D2").Value = 'Barranquilla, Colombia'
$ws.Range("E2").Value = 'Project: Effect of resource availability on women''s preferences for masculinity faces in interaction with hormonal, cognitive, and socio-contextual factors such as actual resource scarcity and exposure to violence: an experimental study using eye-tracking'

$ws.Range("E3").Value = 'COP\$356.040.884 '

$ws.Range("A4").Value = 'IX \href{https://www.unbosque.edu.co/centro-informacion/convocatoria/xiv-convocatoria-interna-de-investigaciones}{Internal Call for Financing Research and Technological Innovation Projects El Bosque University}, 2024'
$ws.Range("B4").Value = 'Jan. 2024 - Jan. 2026'
$ws.Range("C4").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}'
$ws.Range("D4").Value = 'Bogota, Colombia'
$ws.Range("E4").Value = 'Project: Effect of real and simulated resource control on androphilic women''s preferences for masculinity in men''s faces: an experimental study using eye-tracking'

$ws.Range("E5").Value = 'Role: Principal Researcher'

$ws.Range("E6").Value = 'COP\$90.000.000'

$ws.Range("A7").Value = '\href{https://minciencias.gov.co/convocatorias/oportunidades-formacion/convocatoria-programa-estancias-postdoctorales-en-entidades}{Call for Postdoctoral Fellowship Program in SNCTeI entities 2019}'
$ws.Range("B7").Value = 'Jan. 2021 - Jan. 2022'
$ws.Range("C7").Value = '\href{https://minciencias.gov.co/}{Minciencias}'
$ws.Range("D7").Value = 'Barranquilla, Colombia'
$ws.Range("E7").Value = 'Project: Feasibility of new interventions to improve the implementation of sexual and reproductive health programs in Colombia. '

$ws.Range("E8").Value = 'COP\$192.000.000'

$ws.Range("A9").Value = '\href{https://minciencias.gov.co/convocatorias/vocaciones-cientificas-ctei/convocatoria-para-el-fortalecimiento-proyectos-en}{Call for the strengthening of projects in execution of CTeI in health sciences with young talent and regional impact 2020}'
$ws.Range("B9").Value = 'Jan. 2021 - Jan. 2022'
$ws.Range("C9").Value = '\href{https://minciencias.gov.co/}{Minciencias}'
$ws.Range("D9").Value = 'Bogota, Colombia'
$ws.Range("E9").Value = 'Project: Attentional biases and their relationship with heart rate variability as predictors of emotional state in people without affective disorders in the city of Bogotá.'

$ws.Range("E10").Value = 'Role: Principal Researcher'

$ws.Range("E11").Value = 'COP\$76.000.000'

$ws.Range("A12").Value = 'IX \href{https://www.unbosque.edu.co/investigaciones/convocatorias-investigacion}{Internal Call for Financing Research and Technological Innovation Projects El Bosque University}, 2017'
$ws.Range("B12").Value = 'Jan. 2018 - Dic. 2021'
$ws.Range("C12").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}'
$ws.Range("D12").Value = 'Bogota, Colombia'
$ws.Range("E12").Value = 'Project: Perceivable signs of physical and mental health in faces, voices and body odors, and their relationship to hormone levels'

$ws.Range("E13").Value = 'Role: Co-researcher'

$ws.Range("E14").Value = 'COP\$136.586.537'

$ws.Range("A15").Value = 'VII \href{https://www.unbosque.edu.co/investigaciones/convocatorias-investigacion}{Internal Call for Financing Research and Technological Innovation Projects El Bosque University}, 2015'
$ws.Range("B15").Value = 'Jan. 2016 - Dic. 2019'
$ws.Range("C15").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}'
$ws.Range("D15").Value = 'Bogota, Colombia'
$ws.Range("E15").Value = 'Project: Differences in the pattern of eye tracking to sexually preferred stimuli in men convicted of sexual offenses and the general population'

$ws.Range("E16").Value = 'Role: Principal Researcher'

$ws.Range("E17").Value = 'COP\$80.000.000'

$ws.Range("A18").Value = 'Convocatoria Interna de Investigación Financiera de la Universidad de San Buenaventura, 2014'
$ws.Range("B18").Value = 'Jun.2014 - Jun.2015'
$ws.Range("C18").Value = '\href{https://www.usbmed.edu.co/}{Universidad San Buenaventura de Medellín}'
$ws.Range("D18").Value = 'Medellín, Colombia'
$ws.Range("E18").Value = 'Project: Mediating factors of Cognitive Reserve and its relationship with the neuropsychological profile of the older adult in the process of normal aging'

$ws.Range("E19").Value = 'Role: Principal Researcher'

$ws.Range("E20").Value = 'COP\$20.000.000'

# Apply currency-style number format to designated "why" continuation cells
$ws.Range("E3").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E6").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E8").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E11").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E14").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E15").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E17").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'
$ws.Range("E20").NumberFormat = '_-[$$-240A]\ * #,##0.00_-;\-[$$-240A]\ * #,##0.00_-;_-[$$-240A]\ * "-"??_-;_-@_-'

# Row heights for wrapped multi-line rows
$ws.Rows.Item(2).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 43.2
$ws.Rows.Item(15).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 72

# Column widths
$ws.Columns.Item(1).ColumnWidth = 80.6640625
$ws.Columns.Item(2).ColumnWidth = 18.44140625
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 15.5546875
$ws.Columns.Item(5).ColumnWidth = 108.6640625

# Selection
$ws.Range("A18").Select()